# Remove a bunch of audio-file rows (shorter sentences only) from the
# metadata sheet, per commit message "removed a bunch of audio files /
# hopefully, the shorter sentences will avoid NaN loss".
#
# The rows below are identified by the sentence id that prefixes the cell
# text ("N|text|text"); they are the rows deleted upstream. We delete by
# current row number, starting from the bottom so earlier deletions don't
# shift the row numbers still queued for removal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(9, 11, 13, 18, 19, 21, 22, 26, 30, 38, 40, 41, 43, 48, 50, 51, 52, 54, 56, 58, 59, 61)
$sortedDescending = $rowsToDelete | Sort-Object -Descending

foreach ($r in $sortedDescending) {
    $ws.Rows.Item($r).Delete()
}

# Match the author's final selection in the saved view.
$ws.Range("B14").Select()
